# Update recalculated TPM-derived NATMI statistics (columns G:J, M:T)
# for rows 2-13 of the Efnb1-Ephb3 ligand-receptor sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 9.546140333333334
$ws.Cells.Item(2, 8).Value = 28.638421
$ws.Cells.Item(2, 9).Value = 0.587227294878132
$ws.Cells.Item(2, 10).Value = 0.587227294878132
$ws.Cells.Item(2, 13).Value = 0.1709536666666667
$ws.Cells.Item(2, 14).Value = 0.512861
$ws.Cells.Item(2, 15).Value = 0.007882947722998253
$ws.Cells.Item(2, 16).Value = 0.007882947722998253
$ws.Cells.Item(2, 17).Value = 1.631947692497889
$ws.Cells.Item(2, 18).Value = 14.687529232481
$ws.Cells.Item(2, 19).Value = 0.004629082067041994
$ws.Cells.Item(2, 20).Value = 0.004629082067041994
$ws.Cells.Item(3, 7).Value = 9.546140333333334
$ws.Cells.Item(3, 8).Value = 28.638421
$ws.Cells.Item(3, 9).Value = 0.587227294878132
$ws.Cells.Item(3, 10).Value = 0.587227294878132
$ws.Cells.Item(3, 15).Value = 0.7927950496303802
$ws.Cells.Item(3, 16).Value = 0.7927950496303802
$ws.Cells.Item(3, 17).Value = 164.1264279976673
$ws.Cells.Item(3, 18).Value = 1477.137851979006
$ws.Cells.Item(3, 19).Value = 0.4655508923872225
$ws.Cells.Item(3, 20).Value = 0.4655508923872225
$ws.Cells.Item(4, 7).Value = 9.546140333333334
$ws.Cells.Item(4, 8).Value = 28.638421
$ws.Cells.Item(4, 9).Value = 0.587227294878132
$ws.Cells.Item(4, 10).Value = 0.587227294878132
$ws.Cells.Item(4, 13).Value = 4.322599666666666
$ws.Cells.Item(4, 14).Value = 12.967799
$ws.Cells.Item(4, 15).Value = 0.1993220026466216
$ws.Cells.Item(4, 16).Value = 0.1993220026466216
$ws.Cells.Item(4, 17).Value = 41.26414302281989
$ws.Cells.Item(4, 18).Value = 371.377287205379
$ws.Cells.Item(4, 19).Value = 0.1170473204238675
$ws.Cells.Item(4, 20).Value = 0.1170473204238675
$ws.Cells.Item(5, 9).Value = 0.2496684258894083
$ws.Cells.Item(5, 10).Value = 0.2496684258894083
$ws.Cells.Item(5, 13).Value = 0.1709536666666667
$ws.Cells.Item(5, 14).Value = 0.512861
$ws.Cells.Item(5, 15).Value = 0.007882947722998253
$ws.Cells.Item(5, 16).Value = 0.007882947722998253
$ws.Cells.Item(5, 17).Value = 0.6938468546567779
$ws.Cells.Item(5, 18).Value = 6.244621691911001
$ws.Cells.Item(5, 19).Value = 0.001968123149369469
$ws.Cells.Item(5, 20).Value = 0.001968123149369469
$ws.Cells.Item(6, 9).Value = 0.2496684258894083
$ws.Cells.Item(6, 10).Value = 0.2496684258894083
$ws.Cells.Item(6, 15).Value = 0.7927950496303802
$ws.Cells.Item(6, 16).Value = 0.7927950496303802
$ws.Cells.Item(6, 19).Value = 0.1979358920941323
$ws.Cells.Item(6, 20).Value = 0.1979358920941323
$ws.Cells.Item(7, 9).Value = 0.2496684258894083
$ws.Cells.Item(7, 10).Value = 0.2496684258894083
$ws.Cells.Item(7, 13).Value = 4.322599666666666
$ws.Cells.Item(7, 14).Value = 12.967799
$ws.Cells.Item(7, 15).Value = 0.1993220026466216
$ws.Cells.Item(7, 16).Value = 0.1993220026466216
$ws.Cells.Item(7, 17).Value = 17.54406466463878
$ws.Cells.Item(7, 18).Value = 157.896581981749
$ws.Cells.Item(7, 19).Value = 0.04976441064590649
$ws.Cells.Item(7, 20).Value = 0.04976441064590649
$ws.Cells.Item(8, 7).Value = 2.210442
$ws.Cells.Item(8, 8).Value = 6.631326
$ws.Cells.Item(8, 9).Value = 0.1359745227725727
$ws.Cells.Item(8, 10).Value = 0.1359745227725727
$ws.Cells.Item(8, 13).Value = 0.1709536666666667
$ws.Cells.Item(8, 14).Value = 0.512861
$ws.Cells.Item(8, 15).Value = 0.007882947722998253
$ws.Cells.Item(8, 16).Value = 0.007882947722998253
$ws.Cells.Item(8, 17).Value = 0.377883164854
$ws.Cells.Item(8, 18).Value = 3.400948483686
$ws.Cells.Item(8, 19).Value = 0.001071880054675826
$ws.Cells.Item(8, 20).Value = 0.001071880054675826
$ws.Cells.Item(9, 7).Value = 2.210442
$ws.Cells.Item(9, 8).Value = 6.631326
$ws.Cells.Item(9, 9).Value = 0.1359745227725727
$ws.Cells.Item(9, 10).Value = 0.1359745227725727
$ws.Cells.Item(9, 15).Value = 0.7927950496303802
$ws.Cells.Item(9, 16).Value = 0.7927950496303802
$ws.Cells.Item(9, 17).Value = 38.004045309204
$ws.Cells.Item(9, 18).Value = 342.036407782836
$ws.Cells.Item(9, 19).Value = 0.107799928529949
$ws.Cells.Item(9, 20).Value = 0.107799928529949
$ws.Cells.Item(10, 7).Value = 2.210442
$ws.Cells.Item(10, 8).Value = 6.631326
$ws.Cells.Item(10, 9).Value = 0.1359745227725727
$ws.Cells.Item(10, 10).Value = 0.1359745227725727
$ws.Cells.Item(10, 13).Value = 4.322599666666666
$ws.Cells.Item(10, 14).Value = 12.967799
$ws.Cells.Item(10, 15).Value = 0.1993220026466216
$ws.Cells.Item(10, 16).Value = 0.1993220026466216
$ws.Cells.Item(10, 17).Value = 9.554855852385998
$ws.Cells.Item(10, 18).Value = 85.993702671474
$ws.Cells.Item(10, 19).Value = 0.02710271418794784
$ws.Cells.Item(10, 20).Value = 0.02710271418794784
$ws.Cells.Item(11, 7).Value = 0.4410293333333333
$ws.Cells.Item(11, 8).Value = 1.323088
$ws.Cells.Item(11, 9).Value = 0.02712975645988715
$ws.Cells.Item(11, 10).Value = 0.02712975645988715
$ws.Cells.Item(11, 13).Value = 0.1709536666666667
$ws.Cells.Item(11, 14).Value = 0.512861
$ws.Cells.Item(11, 15).Value = 0.007882947722998253
$ws.Cells.Item(11, 16).Value = 0.007882947722998253
$ws.Cells.Item(11, 17).Value = 0.07539558164088889
$ws.Cells.Item(11, 18).Value = 0.678560234768
$ws.Cells.Item(11, 19).Value = 0.0002138624519109645
$ws.Cells.Item(11, 20).Value = 0.0002138624519109646
$ws.Cells.Item(12, 7).Value = 0.4410293333333333
$ws.Cells.Item(12, 8).Value = 1.323088
$ws.Cells.Item(12, 9).Value = 0.02712975645988715
$ws.Cells.Item(12, 10).Value = 0.02712975645988715
$ws.Cells.Item(12, 15).Value = 0.7927950496303802
$ws.Cells.Item(12, 16).Value = 0.7927950496303802
$ws.Cells.Item(12, 17).Value = 7.582600568885332
$ws.Cells.Item(12, 18).Value = 68.243405119968
$ws.Cells.Item(12, 19).Value = 0.02150833661907636
$ws.Cells.Item(12, 20).Value = 0.02150833661907636
$ws.Cells.Item(13, 7).Value = 0.4410293333333333
$ws.Cells.Item(13, 8).Value = 1.323088
$ws.Cells.Item(13, 9).Value = 0.02712975645988715
$ws.Cells.Item(13, 10).Value = 0.02712975645988715
$ws.Cells.Item(13, 13).Value = 4.322599666666666
$ws.Cells.Item(13, 14).Value = 12.967799
$ws.Cells.Item(13, 15).Value = 0.1993220026466216
$ws.Cells.Item(13, 16).Value = 0.1993220026466216
$ws.Cells.Item(13, 17).Value = 1.906393249256889
$ws.Cells.Item(13, 18).Value = 17.157539243312
$ws.Cells.Item(13, 19).Value = 0.005407557388899826
$ws.Cells.Item(13, 20).Value = 0.005407557388899827
